$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# --- Rows with both D (price, text) and E (volume%, text) updates ---
Set-TextValue $ws.Range("D2") '97.623.95'
$ws.Range("E2").Value = '  +0.23%  '
Set-TextValue $ws.Range("D3") '3.471.78'
$ws.Range("E3").Value = '  +4.26%  '
Set-TextValue $ws.Range("D5") '248.83'
$ws.Range("E5").Value = '  +0.45%  '
Set-TextValue $ws.Range("D6") '657.13'
$ws.Range("E6").Value = '  +0.57%  '
Set-TextValue $ws.Range("D7") '1.45'
$ws.Range("E7").Value = '  +5.16%  '
Set-TextValue $ws.Range("D8") '0.418'
$ws.Range("E8").Value = '  -0.04%  '
Set-TextValue $ws.Range("D9") '0.999'
$ws.Range("E9").Value = '  -0.03%  '
Set-TextValue $ws.Range("D10") '1.01'
$ws.Range("E10").Value = '  +1.26%  '
Set-TextValue $ws.Range("D11") '3.469.63'
$ws.Range("E11").Value = '  +4.25%  '
Set-TextValue $ws.Range("D12") '44.14'
$ws.Range("E12").Value = '  +9.77%  '
Set-TextValue $ws.Range("D13") '0.207'
$ws.Range("E13").Value = '  +0.98%  '
Set-TextValue $ws.Range("D14") '97.375.09'
$ws.Range("E14").Value = '  +0.08%  '
Set-TextValue $ws.Range("D15") '6.15'
$ws.Range("E15").Value = '  +2.05%  '
Set-TextValue $ws.Range("D16") '4.121.82'
$ws.Range("E16").Value = '  +4.18%  '
Set-TextValue $ws.Range("D17") '0.0000256'
$ws.Range("E17").Value = '  +1.77%  '
Set-TextValue $ws.Range("D18") '8.71'
$ws.Range("E18").Value = '  +2.48%  '
Set-TextValue $ws.Range("D19") '3.454.51'
$ws.Range("E19").Value = '  +3.68%  '
Set-TextValue $ws.Range("D20") '18.39'
$ws.Range("E20").Value = '  +10.09%  '
Set-TextValue $ws.Range("D21") '11.86'
$ws.Range("E21").Value = '  +12.51%  '
Set-TextValue $ws.Range("D23") '516.98'
$ws.Range("E23").Value = '  +4.38%  '
Set-TextValue $ws.Range("D26") '6.74'
$ws.Range("E26").Value = '  +4.96%  '
Set-TextValue $ws.Range("D27") '96.20'
$ws.Range("E27").Value = '  +3.24%  '
Set-TextValue $ws.Range("D28") '12.44'
$ws.Range("E28").Value = '  +3.88%  '
Set-TextValue $ws.Range("D29") '3.651.25'
$ws.Range("E29").Value = '  +3.95%  '
Set-TextValue $ws.Range("D30") '12.13'
$ws.Range("E30").Value = '  +11.76%  '
Set-TextValue $ws.Range("D35") '0.584'
$ws.Range("E35").Value = '  +7.01%  '
Set-TextValue $ws.Range("D36") '30.68'
$ws.Range("E36").Value = '  +8.86%  '
Set-TextValue $ws.Range("D37") '0.999'
$ws.Range("E37").Value = '  -0.28%  '
Set-TextValue $ws.Range("D38") '7.81'
$ws.Range("E38").Value = '  +3.70%  '
Set-TextValue $ws.Range("D43") '0.900'
$ws.Range("E43").Value = '  +8.34%  '
Set-TextValue $ws.Range("D44") '24.35'
$ws.Range("E44").Value = '  -0.93%  '

# --- Rows with only E (volume%) updated ---
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E24").Value = '  +2.04%  '
$ws.Range("E25").Value = '  +0.07%  '
$ws.Range("E31").Value = '  +14.00%  '
$ws.Range("E32").Value = '  +0.39%  '
$ws.Range("E33").Value = '  -2.84%  '
$ws.Range("E34").Value = '  -1.06%  '
$ws.Range("E39").Value = '  +2.95%  '
$ws.Range("E40").Value = '  +3.70%  '

# --- Rows with only D (price) updated ---
Set-TextValue $ws.Range("D22") '0.502'

# --- Rows 41-51: coin name/link/price/volume reshuffle ---
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range("D41") '1.00'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-TextValue $ws.Range("D42") '515.37'
$ws.Range("E42").Value = '  +2.56%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue $ws.Range("D45") '0.0423'
$ws.Range("E45").Value = '  +3.53%  '
$ws.Range("B46").Value = 'ImmutableX'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D46") '1.71'
$ws.Range("E46").Value = '  +5.13%  '
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws.Range("D47") '3.39'
$ws.Range("E47").Value = '  +8.36%  '
$ws.Range("B48").Value = 'MantraDAO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
Set-TextValue $ws.Range("D48") '3.63'
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue $ws.Range("D49") '5.62'
$ws.Range("E49").Value = '  +3.23%  '
$ws.Range("B50").Value = 'Stacks'
$ws.Range("C50").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Range("D50") '2.21'
$ws.Range("E50").Value = '  +12.20%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue $ws.Range("D51") '8.49'
$ws.Range("E51").Value = '  -0.90%  '

Write-Host "edit complete"
